$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$overview.Range("B3").Value = "Ready for handoff"
$overview.Range("C3").Value = "Ready for handoff"
$overview.Range("D3").Value = "2016-03-24 18:22:09"

$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("E3").Value = "2016-03-24 18:22:02"

$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("E3").Value = "2016-03-24 18:22:09"
